$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 used to hold "password hashing" (unstyled). The task list entry was
# reworded to "authentication" and given the same green highlight used by the
# other task rows above it.
$ws.Range("A18").Value = "authentication"
$ws.Range("A18").Interior.Color = 5296274

# A brand new task row was appended below it for "authorization", left with
# the default (no) fill.
$ws.Range("A19").Value = "authorization"
$ws.Range("A19").Interior.ColorIndex = -4142

# Leave the selection where the author left off editing.
$ws.Range("D14").Select()
